# Daily refresh of the cryptos ranking table (GitHub Actions symbol-list update).
# Updates Price (D) and Volume(1h) (E) figures for most rows, and fixes a few
# Coin/Link (B/C) cells whose rows shifted in the upstream ranking.
# Numeric-looking text (plain numbers and percentages) is written with a
# leading apostrophe so Excel keeps storing it as text instead of coercing it
# to a number, matching the source data's text-formatted columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.13"
$ws.Range("E2").Value = "'0.39%"
$ws.Range("D3").Value = "'29.75"
$ws.Range("E3").Value = "'-0.61%"
$ws.Range("D4").Value = "'5.155"
$ws.Range("E4").Value = "'0.60%"
$ws.Range("D5").Value = "'0.05803"
$ws.Range("E5").Value = "'1.46%"
$ws.Range("D6").Value = "'6.670"
$ws.Range("E6").Value = "'1.69%"
$ws.Range("E7").Value = "'6.93%"
$ws.Range("D8").Value = "'0.8520"
$ws.Range("E8").Value = "'-0.52%"
$ws.Range("D9").Value = "'0.8589"
$ws.Range("E9").Value = "'-1.64%"
$ws.Range("D10").Value = "'0.1374"
$ws.Range("E10").Value = "'2.01%"
$ws.Range("D11").Value = "'0.07099"
$ws.Range("E11").Value = "'2.64%"
$ws.Range("D12").Value = "'0.03203"
$ws.Range("E12").Value = "'10.75%"
$ws.Range("D13").Value = "'0.09378"
$ws.Range("E13").Value = "'0.03%"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("E14").Value = "'1.41%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006029"
$ws.Range("E15").Value = "'1.14%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005927"
$ws.Range("E16").Value = "'-0.91%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.497"
$ws.Range("E17").Value = "'-0.33%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.211"
$ws.Range("E18").Value = "'1.33%"
$ws.Range("D19").Value = "'0.3197"
$ws.Range("E19").Value = "'1.65%"
$ws.Range("D20").Value = "'0.03352"
$ws.Range("E20").Value = "'-0.40%"
$ws.Range("D21").Value = "'0.1296"
$ws.Range("E21").Value = "'-0.48%"
$ws.Range("D22").Value = "'3.485"
$ws.Range("E22").Value = "'-3.29%"
$ws.Range("D23").Value = "'0.04139"
$ws.Range("E23").Value = "'-0.43%"
$ws.Range("D24").Value = "'0.1381"
$ws.Range("E24").Value = "'0.46%"
$ws.Range("D25").Value = "'0.001226"
$ws.Range("E25").Value = "'1.31%"
$ws.Range("D26").Value = "'0.004141"
$ws.Range("E26").Value = "'-7.57%"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'1.91%"
$ws.Range("E28").Value = "'4.20%"
$ws.Range("D40").Value = "'0.03756"
$ws.Range("E40").Value = "'-0.20%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1071"
$ws.Range("E41").Value = "'0.57%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003530"
$ws.Range("E42").Value = "'-38.60%"
$ws.Range("D43").Value = "'0.002430"
$ws.Range("E43").Value = "'9.46%"
$ws.Range("D44").Value = "'0.009559"
$ws.Range("E44").Value = "'0.17%"
$ws.Range("D45").Value = "'0.00005292"
$ws.Range("E45").Value = "'4.48%"
$ws.Range("E46").Value = "'0.31%"
$ws.Range("D47").Value = "'0.05799"
$ws.Range("E47").Value = "'-27.27%"
$ws.Range("E48").Value = "'-20.54%"
$ws.Range("E49").Value = "'0.31%"
$ws.Range("E50").Value = "'0.31%"
